$d = $word.ActiveDocument

# Locate the "react-colorful" list item paragraph - the new "@reduxjs/toolkit"
# entry should be inserted immediately after it, as a new item in the same
# bulleted/numbered list (same style + numbering).
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd() -eq "react-colorful") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    # Duplicate the existing list-item paragraph's range and use
    # InsertParagraphAfter so the new paragraph inherits the same
    # paragraph style (ListParagraph) and numbering (numId 1) as its
    # neighbours, then fill in the new text.
    $dup = $target.Range.Duplicate
    $dup.InsertParagraphAfter()

    $newPara = $target.Next()
    $newPara.Range.Text = "@reduxjs/toolkit "
}
